# Apply the diff:
#  - Insert a new "Player Info" sheet as the first sheet with player
#    metadata (ID, NAME, BATTING_HAND, BOWL_STYLE).
#  - Rename MATCH_CARD_LINK -> MATCH_CODE on "ODI Batting" (col D) and
#    "ODI Bowling" (col B), replacing the full scorecard URL with just
#    the numeric match code.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Player Info" sheet (becomes the first tab) -----------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# text number format forces text storage so "5928" isn't silently coerced
# to a number (every other cell in these sheets is text, incl. numeric-
# looking ones like MATCH_NUMBER "1"); reset the style afterwards so the
# cell doesn't pick up an extra (unstyled) format like the rest of the data
# row
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5928"
$playerInfo.Range("A2").Style = "Normal"
$playerInfo.Range("B2").Value = "Chemar Keron Holder"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Fast"

# match the bold/bordered header style used on the other sheets
# (bold font, thin box border, centered + top-aligned — same look as the
# MATCH_NUMBER-style headers on "ODI Batting" / "ODI Bowling")
$hdr = $playerInfo.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108  # xlCenter
$hdr.VerticalAlignment = -4160    # xlTop
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

# --- 2. ODI Batting: MATCH_CARD_LINK -> MATCH_CODE -------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"
$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4443"
$battingSheet.Range("D2").Style = "Normal"

# --- 3. ODI Bowling: MATCH_CARD_LINK -> MATCH_CODE -------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"
$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4443"
$bowlingSheet.Range("B2").Style = "Normal"
